# Auto-generated: refresh market-data-driven columns (H:N) on each
# tradecraft-leve sheet, per the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Range("H33").Value = 6812.3125
$ws.Range("I33").Value = 56.714287
$ws.Range("J33").Value = 12066.667
$ws.Range("K33").Value = 56.714287
$ws.Range("L33").Value = 12066.667
$ws.Range("M33").Value = 172.285713
$ws.Range("N33").Value = -12524.667

# row 64
$ws.Range("H64").Value = 3774
$ws.Range("I64").Value = 2990
$ws.Range("J64").Value = 4950
$ws.Range("K64").Value = 2990
$ws.Range("L64").Value = 4950
$ws.Range("M64").Value = -2742
$ws.Range("N64").Value = -5446

# row 67
$ws.Range("H67").Value = 3774
$ws.Range("I67").Value = 2990
$ws.Range("J67").Value = 4950
$ws.Range("K67").Value = 2990
$ws.Range("L67").Value = 4950
$ws.Range("M67").Value = -2132
$ws.Range("N67").Value = -6666

# row 127
$ws.Range("H127").Value = 559.35297
$ws.Range("I127").Value = 325.75
$ws.Range("J127").Value = 1120
$ws.Range("K127").Value = 977.25
$ws.Range("L127").Value = 3360
$ws.Range("M127").Value = 3982.75
$ws.Range("N127").Value = -13280

# row 137
$ws.Range("H137").Value = 947.17145
$ws.Range("I137").Value = 821.5862
$ws.Range("J137").Value = 1554.1666
$ws.Range("K137").Value = 2464.7586
$ws.Range("L137").Value = 4662.4998
$ws.Range("M137").Value = 85.24139999999989
$ws.Range("N137").Value = -9762.4998

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 347092.75
$ws.Range("I32").Value = 2442.2856
$ws.Range("J32").Value = 4770107
$ws.Range("K32").Value = 2442.2856
$ws.Range("L32").Value = 4770107
$ws.Range("M32").Value = -2155.2856
$ws.Range("N32").Value = -4770681

# row 45
$ws.Range("H45").Value = 3367.0454
$ws.Range("I45").Value = 3150.7334
$ws.Range("J45").Value = 3830.5715
$ws.Range("K45").Value = 3150.7334
$ws.Range("L45").Value = 3830.5715
$ws.Range("M45").Value = -2773.7334
$ws.Range("N45").Value = -4584.5715

# row 63
$ws.Range("H63").Value = 5322.273
$ws.Range("I63").Value = 3797.8572
$ws.Range("J63").Value = 7990
$ws.Range("K63").Value = 3797.8572
$ws.Range("L63").Value = 7990
$ws.Range("M63").Value = -3111.8572
$ws.Range("N63").Value = -9362

# row 66
$ws.Range("H66").Value = 5322.273
$ws.Range("I66").Value = 3797.8572
$ws.Range("J66").Value = 7990
$ws.Range("K66").Value = 18989.286
$ws.Range("L66").Value = 39950
$ws.Range("M66").Value = -15557.286
$ws.Range("N66").Value = -46814

# row 92
$ws.Range("H92").Value = 48663.332
$ws.Range("J92").Value = 48663.332
$ws.Range("L92").Value = 48663.332
$ws.Range("N92").Value = -53655.332

# row 110
$ws.Range("H110").Value = 758.1724
$ws.Range("I110").Value = 757.4545000000001
$ws.Range("J110").Value = 760.4286
$ws.Range("K110").Value = 757.4545000000001
$ws.Range("L110").Value = 760.4286
$ws.Range("M110").Value = 1287.5455
$ws.Range("N110").Value = -4850.4286

# row 122
$ws.Range("H122").Value = 13615.173
$ws.Range("I122").Value = 15429.934
$ws.Range("K122").Value = 46289.802
$ws.Range("M122").Value = -43839.802

$ws = $wb.Worksheets.Item("BSM")
# row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# row 86
$ws.Range("H86").Value = 27781928
$ws.Range("I86").Value = 40002404
$ws.Range("J86").Value = 8117.364
$ws.Range("K86").Value = 40002404
$ws.Range("L86").Value = 8117.364
$ws.Range("M86").Value = -40001281
$ws.Range("N86").Value = -10363.364

# row 89
$ws.Range("H89").Value = 27781928
$ws.Range("I89").Value = 40002404
$ws.Range("J89").Value = 8117.364
$ws.Range("K89").Value = 200012020
$ws.Range("L89").Value = 40586.82
$ws.Range("M89").Value = -200006404
$ws.Range("N89").Value = -51818.82

# row 107
$ws.Range("H107").Value = 636.67566
$ws.Range("I107").Value = 518.5833
$ws.Range("J107").Value = 4888
$ws.Range("K107").Value = 518.5833
$ws.Range("L107").Value = 4888
$ws.Range("M107").Value = 1401.4167
$ws.Range("N107").Value = -8728

$ws = $wb.Worksheets.Item("CRP")
# row 99
$ws.Range("H99").Value = 951.3333
$ws.Range("I99").Value = 945.25
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 945.25
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 552.75
$ws.Range("N99").Value = -3996

# row 122
$ws.Range("H122").Value = 1530.85
$ws.Range("I122").Value = 1547.2354
$ws.Range("J122").Value = 1438
$ws.Range("K122").Value = 4641.706200000001
$ws.Range("L122").Value = 4314
$ws.Range("M122").Value = -2191.706200000001
$ws.Range("N122").Value = -9214

# row 126
$ws.Range("H126").Value = 951.3333
$ws.Range("I126").Value = 945.25
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 2835.75
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -365.75
$ws.Range("N126").Value = -7940

# row 132
$ws.Range("H132").Value = 28257.73
$ws.Range("I132").Value = 531.3214
$ws.Range("J132").Value = 114517.664
$ws.Range("K132").Value = 1593.9642
$ws.Range("L132").Value = 343552.992
$ws.Range("M132").Value = 936.0357999999999
$ws.Range("N132").Value = -348612.992

$ws = $wb.Worksheets.Item("CUL")
# row 38
$ws.Range("H38").Value = 1788.2858
$ws.Range("I38").Value = 1264
$ws.Range("J38").Value = 2050.4285
$ws.Range("K38").Value = 3792
$ws.Range("L38").Value = 6151.2855
$ws.Range("M38").Value = -3445
$ws.Range("N38").Value = -6845.2855

# row 48
$ws.Range("H48").Value = 1646.2106
$ws.Range("I48").Value = 500
$ws.Range("J48").Value = 1951.8667
$ws.Range("K48").Value = 1500
$ws.Range("L48").Value = 5855.6001
$ws.Range("M48").Value = -1250
$ws.Range("N48").Value = -6355.6001

# row 55
$ws.Range("H55").Value = 2609.0908
$ws.Range("I55").Value = 666.6667
$ws.Range("J55").Value = 3337.5
$ws.Range("K55").Value = 2000.0001
$ws.Range("L55").Value = 10012.5
$ws.Range("M55").Value = -1823.0001
$ws.Range("N55").Value = -10366.5

# row 113
$ws.Range("H113").Value = 1063.96
$ws.Range("I113").Value = 948
$ws.Range("J113").Value = 1113.6571
$ws.Range("K113").Value = 2844
$ws.Range("L113").Value = 3340.9713
$ws.Range("M113").Value = -674
$ws.Range("N113").Value = -7680.971299999999

# row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

# row 97
$ws.Range("H97").Value = 3253.375
$ws.Range("I97").Value = 2604
$ws.Range("K97").Value = 2604
$ws.Range("M97").Value = -2108

# row 102
$ws.Range("H102").Value = 1716.3158
$ws.Range("I102").Value = 1777.8462
$ws.Range("K102").Value = 1777.8462
$ws.Range("M102").Value = -155.8462

# row 126
$ws.Range("H126").Value = 15152339
$ws.Range("I126").Value = 833.3333
$ws.Range("J126").Value = 20834154
$ws.Range("K126").Value = 2499.9999
$ws.Range("L126").Value = 62502462
$ws.Range("M126").Value = -29.9998999999998
$ws.Range("N126").Value = -62507402

$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 1411.9487
$ws.Range("I40").Value = 1392.1786
$ws.Range("K40").Value = 1392.1786
$ws.Range("M40").Value = -1256.1786

# row 41
$ws.Range("H41").Value = 5000001
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

# row 132
$ws.Range("H132").Value = 4609.3403
$ws.Range("I132").Value = 4691.8096
$ws.Range("J132").Value = 3916.6
$ws.Range("K132").Value = 14075.4288
$ws.Range("L132").Value = 11749.8
$ws.Range("M132").Value = -11545.4288
$ws.Range("N132").Value = -16809.8

$ws = $wb.Worksheets.Item("WVR")
# row 42
$ws.Range("H42").Value = 5000
$ws.Range("J42").Value = 5000
$ws.Range("L42").Value = 5000
$ws.Range("N42").Value = -5756

# row 113
$ws.Range("H113").Value = 445.95746
$ws.Range("J113").Value = 671.2857
$ws.Range("L113").Value = 2013.8571
$ws.Range("N113").Value = -6353.8571

# row 126
$ws.Range("H126").Value = 28572622
$ws.Range("I126").Value = 52632668
$ws.Range("J126").Value = 1318.75
$ws.Range("K126").Value = 157898004
$ws.Range("L126").Value = 3956.25
$ws.Range("M126").Value = -157895534
$ws.Range("N126").Value = -8896.25

# row 136
$ws.Range("H136").Value = 31996.906
$ws.Range("I136").Value = 37570.406
$ws.Range("K136").Value = 112711.218
$ws.Range("M136").Value = -110161.218

